$wb = $excel.ActiveWorkbook
$wsTest = $wb.Worksheets.Item("Test")
$wsConfig = $wb.Worksheets.Item("Configuration")

# Add a new "Available Options" column to the Configuration sheet, listing
# the allowed values for each configurable property.
# Set values in the same order the strings were originally authored so the
# shared-strings table lines up (Browser options, Environment options, then
# the column header).
$wsConfig.Range("C2").Value = "Chrome, Firefox"
$wsConfig.Range("C4").Value = "Automation, Staging, UAT, Production"
$wsConfig.Range("C1").Value = "Available Options"
$wsConfig.Range("C3").Value = ""

# Give the new column the same bordered look as the rest of the table.
$wsConfig.Range("A2").Copy()
$wsConfig.Range("C1:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsConfig.Columns.Item(3).ColumnWidth = 36.33

# Restore the selections that were active when the file was saved.
$wsConfig.Range("F6").Select()
$wsTest.Activate()
$wsTest.Range("D10").Select()
